# Revised accomplishment report encoding form
# Append a new row (No. 8 / "sample comments") to the "General Observations
# and Recommendations" table, right after the existing row No. 7
# ("Kindly improve the statement of the GAD mandate under income
# generating projects (page 3").

$d = $word.ActiveDocument

# Locate the table that holds the GAD-mandate comment (row 7) so the
# edit is anchored to content rather than a hard-coded table index.
$targetTable = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $candidate = $d.Tables.Item($i)
    if ($candidate.Range.Text -like "*Kindly improve the statement of the GAD mandate under income generating projects*") {
        $targetTable = $candidate
        break
    }
}
if ($targetTable -eq $null) {
    $targetTable = $d.Tables.Item(1)
}

# Add a brand-new row at the end of the table and fill in its two cells
# to match the existing "No." / comment column layout.
$newRow = $targetTable.Rows.Add()
$newRow.Cells.Item(1).Range.Text = "8"
$newRow.Cells.Item(2).Range.Text = "sample comments"
